$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns appended to the header row (I, J, K)
$ws.Range("I1").Value = "availableDate"
$ws.Range("I2").Value = "火水木金土"

$ws.Range("J1").Value = "startTime"
$ws.Range("K1").Value = "endTime"

# Times stored as real Excel time serials (fraction of a day), formatted h:mm
$ws.Range("J2").Value = 0.4375
$ws.Range("J2").NumberFormat = "h:mm"

$ws.Range("K2").Value = 0.79166666666666663
$ws.Range("K2").NumberFormat = "h:mm"

# Match the author's final selection/active cell
$null = $ws.Range("C5").Select()
